# Insert a new weekly price record for "Zapallo italiano" above the existing
# row 303, shifting all the subsequent records down by one row (matches the
# pattern seen in the XML diff: row 303 becomes new data, and every row that
# followed moves down one position, with the former last row re-appearing as
# the new row 359).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 303 (and everything below it) down by one row.
$ws.Rows.Item(303).Insert()

# Populate the newly-inserted row 303 with the new weekly record.
$ws.Cells.Item(303, 1).Value = 7
$ws.Cells.Item(303, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(303, 3).Value = "Ñuble"
$ws.Cells.Item(303, 4).Value = 45211
$ws.Cells.Item(303, 5).Value = 16
$ws.Cells.Item(303, 6).Value = 100112032
$ws.Cells.Item(303, 7).Value = "Zapallo italiano"
$ws.Cells.Item(303, 8).Value = "Sin especificar"
$ws.Cells.Item(303, 9).Value = "Primera"
$ws.Cells.Item(303, 10).Value = 60
$ws.Cells.Item(303, 11).Value = 16000
$ws.Cells.Item(303, 12).Value = 16000
$ws.Cells.Item(303, 13).Value = 16000
$ws.Cells.Item(303, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(303, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(303, 16).Value = 320
$ws.Cells.Item(303, 17).Value = 50
$ws.Cells.Item(303, 18).Value = "Hortaliza"
